$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End | Animal Glue
$ws.Range("H5").Value = 104
$ws.Range("I5").Value = 98.666664
$ws.Range("K5").Value = 98.666664
$ws.Range("M5").Value = 16.333336

# Row 18: You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 1812.6666
$ws.Range("I18").Value = 1812.6666
$ws.Range("K18").Value = 1812.6666
$ws.Range("M18").Value = -1528.6666

# Row 31: Hush Little Wailer | Weak Silencing Potion
$ws.Range("H31").Value = 5089
$ws.Range("I31").Value = 5089
$ws.Range("K31").Value = 15267
$ws.Range("M31").Value = -15037

# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 8781.125
$ws.Range("I43").Value = 7187.375
$ws.Range("K43").Value = 7187.375
$ws.Range("M43").Value = -7118.375

# Row 44: Alive and Unwell | Budding Oak Wand
$ws.Range("H44").Value = 6000
$ws.Range("J44").Value = 6000
$ws.Range("L44").Value = 6000
$ws.Range("N44").Value = -6924

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 4134.6665
$ws.Range("I76").Value = 4134.6665
$ws.Range("K76").Value = 4134.6665
$ws.Range("M76").Value = -3819.6665

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 4134.6665
$ws.Range("I79").Value = 4134.6665
$ws.Range("K79").Value = 4134.6665
$ws.Range("M79").Value = -3042.6665

# Row 96: Scroll Down | Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 1032.1333
$ws.Range("I96").Value = 1032.1333
$ws.Range("K96").Value = 3096.3999
$ws.Range("M96").Value = -1723.3999

# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = ""

# Row 114: Conserving Combat | Bluespirit Codex
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1170.4445
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 947.375
$ws.Range("I135").Value = 566.43475
$ws.Range("K135").Value = 5097.91275
$ws.Range("M135").Value = -2562.91275

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 4423.375
$ws.Range("I32").Value = 3022.7754
$ws.Range("K32").Value = 3022.7754
$ws.Range("M32").Value = -2735.7754

# Row 56: Feasting the Night Away | Hells' Kitchen
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2048.1875
$ws.Range("I61").Value = 2004.2858
$ws.Range("K61").Value = 2004.2858
$ws.Range("M61").Value = -1792.2858

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2048.1875
$ws.Range("I136").Value = 2004.2858
$ws.Range("K136").Value = 6012.857400000001
$ws.Range("M136").Value = -3462.857400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2642.7334
$ws.Range("I134").Value = 1785.5454
$ws.Range("K134").Value = 5356.6362
$ws.Range("M134").Value = -2821.6362

# Row 140: Ceremonial Teeth | Ra'Kaznar Twinfangs
$ws.Range("H140").Value = 140853.33
$ws.Range("J140").Value = 140853.33
$ws.Range("L140").Value = 140853.33
$ws.Range("N140").Value = -151213.33

$ws = $wb.Worksheets.Item("CRP")
# Row 41: The Lone Bowman | Oak Longbow
$ws.Range("H41").Value = 21246.666
$ws.Range("I41").Value = 4900
$ws.Range("K41").Value = 4900
$ws.Range("M41").Value = -4472

# Row 47: Grippy When Wet | Mythril Cavalry Bow
$ws.Range("H47").Value = 19035
$ws.Range("J47").Value = 19035
$ws.Range("L47").Value = 19035
$ws.Range("N47").Value = -20167

# Row 51: Greenstone for Greenhorns | Jade Crook
$ws.Range("H51").Value = 29965.666
$ws.Range("J51").Value = 29965.666
$ws.Range("L51").Value = 29965.666
$ws.Range("N51").Value = -31437.666

# Row 60: Bowing to Greater Power | Yew Longbow
$ws.Range("H60").Value = 14530.462
$ws.Range("I60").Value = 10435.556
$ws.Range("J60").Value = 23744
$ws.Range("K60").Value = 10435.556
$ws.Range("L60").Value = 23744
$ws.Range("M60").Value = -9924.556
$ws.Range("N60").Value = -24766

# Row 61: Incant Now, Think Later | Jade Crook
$ws.Range("H61").Value = 29965.666
$ws.Range("J61").Value = 29965.666
$ws.Range("L61").Value = 29965.666
$ws.Range("N61").Value = -30661.666

# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 201499.5
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 201499.5
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

# Row 69: Landing the Big One | Cedar Fishing Rod
$ws.Range("H69").Value = 22993.334
$ws.Range("I69").Value = 21990
$ws.Range("K69").Value = 21990
$ws.Range("M69").Value = -21241

# Row 72: Fishing for Profits (L) | Cedar Fishing Rod
$ws.Range("H72").Value = 22993.334
$ws.Range("I72").Value = 21990
$ws.Range("K72").Value = 65970
$ws.Range("M72").Value = -62226

$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 989.913
$ws.Range("J131").Value = 989.913
$ws.Range("L131").Value = 2969.739
$ws.Range("N131").Value = -13049.739

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success | Copper Ring
$ws.Range("H11").Value = 7667000
$ws.Range("J11").Value = 13000000
$ws.Range("L11").Value = 13000000
$ws.Range("N11").Value = -13000278

# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 6000
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 6000
$ws.Range("N70").Value = -6540

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 6000
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 6000
$ws.Range("N73").Value = -7872

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 60916.824
$ws.Range("I122").Value = 1424.4546
$ws.Range("J122").Value = 169986.17
$ws.Range("K122").Value = 4273.3638
$ws.Range("L122").Value = 509958.51
$ws.Range("M122").Value = -1823.3638
$ws.Range("N122").Value = -514858.51

$ws = $wb.Worksheets.Item("LTW")
# Row 20: Choke Hold | Hard Leather Choker
$ws.Range("H20").Value = 29566.5
$ws.Range("J20").Value = 30006
$ws.Range("L20").Value = 30006
$ws.Range("N20").Value = -30458

# Row 41: The Hand that Bleeds | Fingerless Boarskin Gloves
$ws.Range("H41").Value = 24000
$ws.Range("I41").Value = 24000
$ws.Range("K41").Value = 24000
$ws.Range("M41").Value = -23562

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 1359.6154
$ws.Range("I93").Value = 1153.2222
$ws.Range("K93").Value = 1153.2222
$ws.Range("M93").Value = 94.77780000000007

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3548.3333
$ws.Range("I132").Value = 2895
$ws.Range("J132").Value = 3875
$ws.Range("K132").Value = 8685
$ws.Range("L132").Value = 11625
$ws.Range("M132").Value = -6155
$ws.Range("N132").Value = -16685

$ws = $wb.Worksheets.Item("WVR")
# Row 48: In over Your Head | Linen Cowl
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 50000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 50000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -49431
$ws.Range("N48").Value = ""

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2918.2
$ws.Range("I132").Value = 2304
$ws.Range("K132").Value = 6912
$ws.Range("M132").Value = -4382
